$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The BOM resistor rows (17-19) are being re-split: the old combined
# "R511, R512, R514, R515, R518, R519, R520 / 16k" line is split into a
# new "R518, R519, R520 / 270k" part plus a reduced "R511, R512, R514,
# R515 / 16k" line, and the remaining rows shift down, adding a new
# row 20 for the part that used to be on row 18.

# --- Row 17: new part "R518, R519, R520" / R 270k 0402 ---
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "R518, R519, R520"
$ws.Range("C17").Value = "Resistor"
$ws.Range("D17").Value = "R 270k 0402"
$ws.Range("E17").Value = "RES 270K OHM 1% 1/16W 0402"
$ws.Range("F17").Value = "YAGEO"
$ws.Range("G17").Value = "RC0402FR-07270KL"
$ws.Range("H17").Value = "R_270k_0402"
$ws.Range("I17").Value = "Digi-Key"
$ws.Range("J17").Value = "311-270KLRCT-ND"

# --- Row 18: "R500, R501, R503, R504" / R 10k 0402 (Manufacturer fixed to YAGEO) ---
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "R500, R501, R503, R504"
$ws.Range("C18").Value = "Resistor"
$ws.Range("D18").Value = "R 10k 0402"
$ws.Range("E18").Value = "RES 10K OHM 1% 1/16W 0402"
$ws.Range("F18").Value = "YAGEO"
$ws.Range("G18").Value = "RC0402FR-0710KL"
$ws.Range("H18").Value = "R_10k_0402"
$ws.Range("I18").Value = "Digi-Key"
$ws.Range("J18").Value = "311-10.0KLRCT-ND"

# --- Row 19: "R511, R512, R514, R515" / R 16k 0402 (designators R518-520 split off) ---
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "R511, R512, R514, R515"
$ws.Range("C19").Value = "Resistor"
$ws.Range("D19").Value = "R 16k 0402"
$ws.Range("E19").Value = "RES 16K OHM 1% 1/16W 0402"
$ws.Range("F19").Value = "Stackpole Electronics Inc"
$ws.Range("G19").Value = "RMCF0402FT16K0"
$ws.Range("H19").Value = "R_16k_0402"
$ws.Range("I19").Value = "Digi-Key"
$ws.Range("J19").Value = "RMCF0402FT16K0CT-ND"

# --- Row 20 (new row): "R505, R507, R508, R516, R517" / R 0R0 0402 ---
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "R505, R507, R508, R516, R517"
$ws.Range("C20").Value = "Resistor"
$ws.Range("D20").Value = "R 0R0 0402"
$ws.Range("E20").Value = "RES 0 OHM JUMPER 1/16W 0402"
$ws.Range("F20").Value = "Stackpole Electronics Inc"
$ws.Range("G20").Value = "RMCF0402ZT0R00"
$ws.Range("H20").Value = "R_0R0_0402"
$ws.Range("I20").Value = "Digi-Key"
$ws.Range("J20").Value = "RMCF0402ZT0R00CT-ND"

# Re-apply the standard BOM-row formatting (quantity style in col A,
# text style in cols B:J) to every rewritten/new row, since writing
# .Value resets direct formatting. Row 16 is untouched and keeps the
# canonical per-column styles, so use it as the format source.
$ws.Range("A16:J16").Copy()
$ws.Range("A17:J20").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
